$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" column (M) mirroring the formatting of the existing
# "2020" column (L), row by row, then fill in the new values.

# Row 2: blank separator cell, just carry the format/border over.
$ws.Range("L2").Copy()
$ws.Range("M2").PasteSpecial(-4122)

# Row 3: year header.
$ws.Range("L3").Copy()
$ws.Range("M3").PasteSpecial(-4122)
$ws.Range("M3").Value = 2021

# Rows 4-10: data values.
$ws.Range("L4").Copy()
$ws.Range("M4").PasteSpecial(-4122)
$ws.Range("M4").Value = 952

$ws.Range("L5").Copy()
$ws.Range("M5").PasteSpecial(-4122)
$ws.Range("M5").Value = 10437

$ws.Range("L6").Copy()
$ws.Range("M6").PasteSpecial(-4122)
$ws.Range("M6").Value = 2253

$ws.Range("L7").Copy()
$ws.Range("M7").PasteSpecial(-4122)
$ws.Range("M7").Value = 8184

$ws.Range("L8").Copy()
$ws.Range("M8").PasteSpecial(-4122)
$ws.Range("M8").Value = 14020

$ws.Range("L9").Copy()
$ws.Range("M9").PasteSpecial(-4122)
$ws.Range("M9").Value = 5139

$ws.Range("L10").Copy()
$ws.Range("M10").PasteSpecial(-4122)
$ws.Range("M10").Value = 8881

# Update the active cell / selection to P8.
$ws.Range("P8").Select()
$excel.CutCopyMode = $false
